$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.2722696666666667
$ws.Cells.Item(2, 8).Value = 0.816809
$ws.Cells.Item(2, 9).Value = 0.01577089893809228
$ws.Cells.Item(2, 10).Value = 0.01577089893809228
$ws.Cells.Item(2, 13).Value = 61.04160633333334
$ws.Cells.Item(2, 14).Value = 183.124819
$ws.Cells.Item(2, 15).Value = 0.2043613460574534
$ws.Cells.Item(2, 16).Value = 0.2043613460574534
$ws.Cells.Item(2, 17).Value = 16.61977780917456
$ws.Cells.Item(2, 18).Value = 149.578000282571
$ws.Cells.Item(2, 19).Value = 0.003222962135524601
$ws.Cells.Item(2, 20).Value = 0.0032229621355246

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.2722696666666667
$ws.Cells.Item(3, 8).Value = 0.816809
$ws.Cells.Item(3, 9).Value = 0.01577089893809228
$ws.Cells.Item(3, 10).Value = 0.01577089893809228
$ws.Cells.Item(3, 15).Value = 0.3559304658284363
$ws.Cells.Item(3, 16).Value = 0.3559304658284363
$ws.Cells.Item(3, 17).Value = 28.94620421966467
$ws.Cells.Item(3, 18).Value = 260.515837976982
$ws.Cells.Item(3, 19).Value = 0.005613343405568376
$ws.Cells.Item(3, 20).Value = 0.005613343405568376

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.2722696666666667
$ws.Cells.Item(4, 8).Value = 0.816809
$ws.Cells.Item(4, 9).Value = 0.01577089893809228
$ws.Cells.Item(4, 10).Value = 0.01577089893809228
$ws.Cells.Item(4, 13).Value = 131.3384093333333
$ws.Cells.Item(4, 14).Value = 394.015228
$ws.Cells.Item(4, 15).Value = 0.4397081881141102
$ws.Cells.Item(4, 16).Value = 0.4397081881141103
$ws.Cells.Item(4, 17).Value = 35.75946492971688
$ws.Cells.Item(4, 18).Value = 321.835184367452
$ws.Cells.Item(4, 19).Value = 0.006934593396999301
$ws.Cells.Item(4, 20).Value = 0.006934593396999301

$ws.Cells.Item(5, 9).Value = 0.8050543166133334
$ws.Cells.Item(5, 10).Value = 0.8050543166133333
$ws.Cells.Item(5, 13).Value = 61.04160633333334
$ws.Cells.Item(5, 14).Value = 183.124819
$ws.Cells.Item(5, 15).Value = 0.2043613460574534
$ws.Cells.Item(5, 16).Value = 0.2043613460574534
$ws.Cells.Item(5, 17).Value = 848.3868877070462
$ws.Cells.Item(5, 18).Value = 7635.481989363416
$ws.Cells.Item(5, 19).Value = 0.1645219837924641
$ws.Cells.Item(5, 20).Value = 0.1645219837924641

$ws.Cells.Item(6, 9).Value = 0.8050543166133334
$ws.Cells.Item(6, 10).Value = 0.8050543166133333
$ws.Cells.Item(6, 15).Value = 0.3559304658284363
$ws.Cells.Item(6, 16).Value = 0.3559304658284363
$ws.Cells.Item(6, 19).Value = 0.2865433579293772
$ws.Cells.Item(6, 20).Value = 0.2865433579293772

$ws.Cells.Item(7, 9).Value = 0.8050543166133334
$ws.Cells.Item(7, 10).Value = 0.8050543166133333
$ws.Cells.Item(7, 13).Value = 131.3384093333333
$ws.Cells.Item(7, 14).Value = 394.015228
$ws.Cells.Item(7, 15).Value = 0.4397081881141102
$ws.Cells.Item(7, 16).Value = 0.4397081881141103
$ws.Cells.Item(7, 17).Value = 1825.407144796152
$ws.Cells.Item(7, 18).Value = 16428.66430316537
$ws.Cells.Item(7, 19).Value = 0.3539889748914921
$ws.Cells.Item(7, 20).Value = 0.3539889748914921

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 3.093283333333334
$ws.Cells.Item(8, 8).Value = 9.279850000000001
$ws.Cells.Item(8, 9).Value = 0.1791747844485745
$ws.Cells.Item(8, 10).Value = 0.1791747844485744
$ws.Cells.Item(8, 13).Value = 61.04160633333334
$ws.Cells.Item(8, 14).Value = 183.124819
$ws.Cells.Item(8, 15).Value = 0.2043613460574534
$ws.Cells.Item(8, 16).Value = 0.2043613460574534
$ws.Cells.Item(8, 17).Value = 188.8189835107945
$ws.Cells.Item(8, 18).Value = 1699.37085159715
$ws.Cells.Item(8, 19).Value = 0.03661640012946475
$ws.Cells.Item(8, 20).Value = 0.03661640012946474

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 3.093283333333334
$ws.Cells.Item(9, 8).Value = 9.279850000000001
$ws.Cells.Item(9, 9).Value = 0.1791747844485745
$ws.Cells.Item(9, 10).Value = 0.1791747844485744
$ws.Cells.Item(9, 15).Value = 0.3559304658284363
$ws.Cells.Item(9, 16).Value = 0.3559304658284363
$ws.Cells.Item(9, 17).Value = 328.8607657700334
$ws.Cells.Item(9, 18).Value = 2959.746891930301
$ws.Cells.Item(9, 19).Value = 0.06377376449349077
$ws.Cells.Item(9, 20).Value = 0.06377376449349077

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 3.093283333333334
$ws.Cells.Item(10, 8).Value = 9.279850000000001
$ws.Cells.Item(10, 9).Value = 0.1791747844485745
$ws.Cells.Item(10, 10).Value = 0.1791747844485744
$ws.Cells.Item(10, 13).Value = 131.3384093333333
$ws.Cells.Item(10, 14).Value = 394.015228
$ws.Cells.Item(10, 15).Value = 0.4397081881141102
$ws.Cells.Item(10, 16).Value = 0.4397081881141103
$ws.Cells.Item(10, 17).Value = 406.2669126173111
$ws.Cells.Item(10, 18).Value = 3656.4022135558
$ws.Cells.Item(10, 19).Value = 0.07878461982561893
$ws.Cells.Item(10, 20).Value = 0.07878461982561892
